$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.858.60"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "2.749.90"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("E9").Value = "  -3.19%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -14.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.382"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").Value = "3.236.80"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.86%  "
$ws.Range("D15").Value = "63.510.22"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").Value = "2.752.84"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "355.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.533"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.170"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").Value = "0.0₃0907"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("E29").Value = "  -4.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "168.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.987"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.01%  "
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "330.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0586"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0253"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.51%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.625"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
